$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 157 (shifts the existing rows 157-159 down to 158-160,
# matching the structure in the updated workbook).
$ws.Rows(157).Insert()

# Populate the newly inserted row 157 with the new weekly price record.
$ws.Cells.Item(157, 1).Value = 9
$ws.Cells.Item(157, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(157, 3).Value = "Metropolitana"
$ws.Cells.Item(157, 4).Value = 45041
$ws.Cells.Item(157, 5).Value = 13
$ws.Cells.Item(157, 6).Value = "Fruta"
$ws.Cells.Item(157, 7).Value = 100101
$ws.Cells.Item(157, 8).Value = "Berries"
$ws.Cells.Item(157, 9).Value = 100101004
$ws.Cells.Item(157, 10).Value = "Frambuesa"
$ws.Cells.Item(157, 11).Value = "Sin especificar"
$ws.Cells.Item(157, 12).Value = "Primera"
$ws.Cells.Item(157, 13).Value = 470
$ws.Cells.Item(157, 14).Value = 7000
$ws.Cells.Item(157, 15).Value = 7500
$ws.Cells.Item(157, 16).Value = 7266
$ws.Cells.Item(157, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(157, 18).Value = "Provincia de Colchagua"
$ws.Cells.Item(157, 19).Value = 3633
$ws.Cells.Item(157, 20).Value = 2

# Make sure the date cell keeps the expected date/time number format (style
# index 2 in the original styles.xml), same as the other "Fecha" cells.
$ws.Cells.Item(157, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
